$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repulled data
$ws.Range("F5").Value = -5
$ws.Range("F6").Value = -10
$ws.Range("F10").Value = -3
$ws.Range("F11").Value = -4
$ws.Range("F16").Value = 2
$ws.Range("F17").Value = -14

$wb.Save()
